$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the contents of row 17 and row 18 for columns
# A (Id), K (Alder-Stadium), P (Lokalnamn), Q (Ost), R (Nord), S (Noggrannhet),
# while column B (Taxonsorteringsordning) is updated to the same new value
# (55684) on both rows.

# Capture current ("before") values for row 17 so we can move them to row 18.
$A17 = $ws.Range("A17").Value()
$K17 = $ws.Range("K17").Value()
$P17 = $ws.Range("P17").Value()
$Q17 = $ws.Range("Q17").Value()
$R17 = $ws.Range("R17").Value()
$S17 = $ws.Range("S17").Value()

# Capture current ("before") values for row 18 so we can move them to row 17.
$A18 = $ws.Range("A18").Value()
$K18 = $ws.Range("K18").Value()
$P18 = $ws.Range("P18").Value()
$Q18 = $ws.Range("Q18").Value()
$R18 = $ws.Range("R18").Value()
$S18 = $ws.Range("S18").Value()

# Row 17 gets row 18's former values.
$ws.Range("A17").Value = $A18
$ws.Range("K17").Value = $K18
$ws.Range("P17").Value = $P18
$ws.Range("Q17").Value = $Q18
$ws.Range("R17").Value = $R18
$ws.Range("S17").Value = $S18

# Row 18 gets row 17's former values.
$ws.Range("A18").Value = $A17
$ws.Range("K18").Value = $K17
$ws.Range("P18").Value = $P17
$ws.Range("Q18").Value = $Q17
$ws.Range("R18").Value = $R17
$ws.Range("S18").Value = $S17

# Column B becomes 55684 on both rows.
$ws.Range("B17").Value = 55684
$ws.Range("B18").Value = 55684
